$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '63.810.38'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.640.44'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '607.87'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.34'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.384'
$ws.Range('E10').Value = '  +6.94%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.41'
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.114.10'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '63.622.85'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.646.85'
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.75'
$ws.Range('E18').Value = '  +2.35%  '
$ws.Range('E19').Value = '  +3.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '346.72'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.90'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.45'
$ws.Range('E24').Value = '  -1.81%  '
$ws.Range('E25').Value = '  +6.94%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.69'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.23'
$ws.Range('E27').Value = '  +7.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '562.32'
$ws.Range('E28').Value = '  +3.83%  '
$ws.Range('E29').Value = '  +2.64%  '
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0₃0854'
$ws.Range('E33').Value = '  +5.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.76'
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.29'
$ws.Range('E35').Value = '  +3.48%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '169.35'
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.405'
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('E39').Value = '  +4.74%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '19.15'
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '165.01'
$ws.Range('E42').Value = '  -6.84%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '40.09'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.80'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.88'
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0566'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('E47').Value = '  -1.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.01'
$ws.Range('E48').Value = '  +14.22%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0245'
$ws.Range('E49').Value = '  +1.37%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '18.80'
$ws.Range('E51').Value = '  -0.83%  '
